{"js": "// The author byline paragraph currently reads \"Student Names\" (built from\n// two runs: \"Student Name\" + \"s\"). Replace it with the actual author names,\n// \"Lachlan Dietrich | Rokhan Khattak\", per the commit's authorship update.\n\nconst body = context.document.body;\nconst results = body.search(\"Student Names\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items,text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find \"Student Names\" text to replace.');\n}\n\n// Replace the found range's text in place (keeps it in its own paragraph).\nresults.items[0].insertText(\"Lachlan Dietrich | Rokhan Khattak\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# The author byline paragraph currently reads \"Student Names\" (built from\n# two runs: \"Student Name\" + \"s\"). Replace it with the actual author names,\n# \"Lachlan Dietrich | Rokhan Khattak\", per the commit's authorship update.\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n$found = $find.Execute(\n    \"Student Names\",   # FindText\n    $false,            # MatchCase\n    $false,            # MatchWholeWord\n    $false,            # MatchWildcards\n    $false,            # MatchSoundsLike\n    $false,            # MatchAllWordForms\n    $true,             # Forward\n    1,                 # Wrap (wdFindContinue)\n    $false,            # Format\n    \"Lachlan Dietrich | Rokhan Khattak\",  # ReplaceWith\n    2                  # Replace (wdReplaceAll)\n)\n\nif (-not $found) {\n    throw 'Could not find \"Student Names\" text to replace.'\n}\n"}
